$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(65).Insert()

$ws.Range("A65").Value = 7
$ws.Range("B65").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C65").Value = "Ñuble"
$ws.Range("D65").Value = 44579
$ws.Range("E65").Value = 16
$ws.Range("F65").Value = 100112017
$ws.Range("G65").Value = "Apio"
$ws.Range("H65").Value = "Americana (o)"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 60
$ws.Range("K65").Value = 8000
$ws.Range("L65").Value = 8500
$ws.Range("M65").Value = 8250
$ws.Range("N65").Value = "$/docena de matas"
$ws.Range("O65").Value = "Provincia del Elquí"
$ws.Range("P65").Value = 1375
$ws.Range("Q65").Value = 6
$ws.Range("R65").Value = "Hortaliza"
